$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated F/G column values for existing rows (per diff)
$updates = @{
    293 = @{ F = 81778 }
    294 = @{ F = 92366 }
    299 = @{ F = 64635; G = 6828 }
    300 = @{ F = 71573 }
    324 = @{ F = 234655; G = 2708 }
    325 = @{ F = 763106; G = 6496 }
    326 = @{ F = 433179; G = 3846 }
    327 = @{ F = 237460; G = 2894 }
    331 = @{ F = 150869; G = 2628 }
    332 = @{ F = 435939; G = 4303 }
    333 = @{ F = 265236; G = 2858 }
    334 = @{ F = 203736; G = 3396 }
    335 = @{ F = 128742; G = 2883 }
    336 = @{ F = 100682; G = 3198 }
    337 = @{ F = 102395 }
    338 = @{ F = 218027; G = 3065 }
    339 = @{ F = 642727; G = 5498 }
    340 = @{ F = 381224; G = 3279 }
    341 = @{ F = 295825; G = 3664 }
    342 = @{ F = 173210; G = 2936 }
    343 = @{ F = 127110; G = 2829 }
    344 = @{ F = 130374; G = 2404 }
    345 = @{ F = 275770; G = 3147 }
    346 = @{ F = 637374; G = 4530 }
    347 = @{ F = 324137; G = 2734 }
    348 = @{ F = 221216; G = 3010 }
    349 = @{ F = 158544; G = 2678 }
    350 = @{ F = 116615; G = 2595 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# Append new row 351 with the next day's data
$ws.Range("A351").Value = 44245
$ws.Range("B351").Value = 287752
$ws.Range("C351").Value = 10283
$ws.Range("D351").Value = 2333
$ws.Range("E351").Value = 6350
$ws.Range("F351").Value = 113492
$ws.Range("G351").Value = 2214

# Match date format style of column A used in prior rows
$ws.Range("A351").NumberFormat = $ws.Range("A350").NumberFormat
